# Fix the product name: insert a hyphen after "245" so it reads
# "245-MS-EI-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME" instead of
# "245MS-EI-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME".
$wb = $excel.ActiveWorkbook

$newName = "245-MS-EI-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-ONTIME"

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$wsInput.Range("B1").Value = $newName
$wsOutput.Range("B1").Value = $newName

# Update selection on the input sheet, then move to / activate the
# output sheet and select its B1 cell, leaving it as the active sheet.
$wsInput.Range("B1").Select()
$wsOutput.Activate()
$wsOutput.Range("B1").Select()
